$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.937.67"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "2.364.91"
$ws.Range("E3").Value = "  -1.67%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.79"
$ws.Range("E5").Value = "  -5.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.46"
$ws.Range("E6").Value = "  +6.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.637"
$ws.Range("E7").Value = "  -1.25%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.626"
$ws.Range("E9").Value = "  -2.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.40"
$ws.Range("E10").Value = "  +1.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0929"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.66"
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("E13").Value = "  -2.33%  "
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.02"
$ws.Range("E15").Value = "  -5.42%  "
$ws.Range("D16").Value = "2.723.31"
$ws.Range("E16").Value = "  -1.63%  "
$ws.Range("D17").Value = "2.364.46"
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("D18").Value = "42.934.54"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("E19").Value = "  -1.02%  "
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "76.43"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.64"
$ws.Range("E22").Value = "  -5.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "267.75"
$ws.Range("E23").Value = "  -0.87%  "
$ws.Range("E24").Value = "  -1.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.47"
$ws.Range("E25").Value = "  -8.83%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.50"
$ws.Range("E27").Value = "  -2.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.52"
$ws.Range("E28").Value = "  -3.32%  "
$ws.Range("E29").Value = "  +2.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.92"
$ws.Range("E30").Value = "  +2.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "168.24"
$ws.Range("E31").Value = "  -3.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0912"
$ws.Range("E32").Value = "  -1.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.03"
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.91"
$ws.Range("E34").Value = "  -6.81%  "
$ws.Range("E35").Value = "  -2.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.120"
$ws.Range("E36").Value = "  +11.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.78"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.88"
$ws.Range("E39").Value = "  -1.62%  "
$ws.Range("E40").Value = "  -5.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.10"
$ws.Range("E41").Value = "  +13.42%  "
$ws.Range("E42").Value = "  -2.49%  "
$ws.Range("E43").Value = "  +2.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.49"
$ws.Range("E44").Value = "  +2.37%  "
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.49"
$ws.Range("E46").Value = "  +4.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "113.85"
$ws.Range("E47").Value = "  -3.66%  "
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "77.32"
$ws.Range("E50").Value = "  +9.70%  "
$ws.Range("E51").Value = "  +0.99%  "
